# Auto-generated update of betting-odds values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.82
$ws.Range("I2").Value = 2.25
$ws.Range("J2").Value = 3.3
$ws.Range("L2").Value = 2.8
$ws.Range("N2").Value = 8.5
$ws.Range("S2").Value = 1.34
$ws.Range("T2").Value = 3
$ws.Range("V2").Value = 2.27
$ws.Range("W2").Value = 11.5
$ws.Range("X2").Value = 16.5
$ws.Range("Y2").Value = 10.25
$ws.Range("Z2").Value = 35
$ws.Range("AA2").Value = 21
$ws.Range("AB2").Value = 25
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 7
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 9.75
$ws.Range("AI2").Value = 12.5
$ws.Range("AJ2").Value = 9
$ws.Range("AK2").Value = 23
$ws.Range("AL2").Value = 16.5
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 4.9
$ws.Range("AO2").Value = 14.5
$ws.Range("AP2").Value = 20
$ws.Range("AQ2").Value = 60
$ws.Range("AR2").Value = 80
$ws.Range("AS2").Value = 200
$ws.Range("AT2").Value = 3
$ws.Range("AU2").Value = 6.7
$ws.Range("AW2").Value = 4.35
$ws.Range("AX2").Value = 11.5
$ws.Range("AY2").Value = 18
$ws.Range("AZ2").Value = 45
$ws.Range("BA2").Value = 70
$ws.Range("BB2").Value = 200

# Row 3
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.85
$ws.Range("I3").Value = 3.1
$ws.Range("J3").Value = 2.47
$ws.Range("K3").Value = 2.42
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 9.75
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = 1.44
$ws.Range("R3").Value = 2.62
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 3.6
$ws.Range("U3").Value = 1.42
$ws.Range("V3").Value = 2.67
$ws.Range("Z3").Value = 21
$ws.Range("AC3").Value = 9.75
$ws.Range("AD3").Value = 8.25
$ws.Range("AE3").Value = 11
$ws.Range("AF3").Value = 32
$ws.Range("AH3").Value = 16
$ws.Range("AI3").Value = 22
$ws.Range("AJ3").Value = 11.25
$ws.Range("AM3").Value = 22
$ws.Range("AN3").Value = 4.5
$ws.Range("AO3").Value = 9.75
$ws.Range("AR3").Value = 45
$ws.Range("AT3").Value = 3.6
$ws.Range("AW3").Value = 5.7
$ws.Range("AY3").Value = 17.5
$ws.Range("AZ3").Value = 65
$ws.Range("BA3").Value = 70
$ws.Range("BC3").Value = 450

# Row 4
$ws.Range("G4").Value = 2.57
$ws.Range("I4").Value = 2.77
$ws.Range("J4").Value = 3.15
$ws.Range("L4").Value = 3.35
$ws.Range("W4").Value = 7
$ws.Range("X4").Value = 11.75
$ws.Range("Y4").Value = 10
$ws.Range("Z4").Value = 28
$ws.Range("AA4").Value = 24
$ws.Range("AB4").Value = 37
$ws.Range("AH4").Value = 7.4
$ws.Range("AI4").Value = 13
$ws.Range("AJ4").Value = 10.25
$ws.Range("AK4").Value = 32
$ws.Range("AL4").Value = 26
$ws.Range("AM4").Value = 40
$ws.Range("AN4").Value = 4.4
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 22
$ws.Range("AQ4").Value = 60
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 300
$ws.Range("AW4").Value = 4.65
$ws.Range("AX4").Value = 15
$ws.Range("AY4").Value = 23
$ws.Range("AZ4").Value = 70
$ws.Range("BA4").Value = 110

# Row 8
$ws.Range("G8").Value = 2.2
$ws.Range("I8").Value = 3.4
$ws.Range("L8").Value = 3.6
$ws.Range("Q8").Value = 1.8
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.2
$ws.Range("W8").Value = 9
$ws.Range("Z8").Value = 21
$ws.Range("AD8").Value = 6.5
$ws.Range("AI8").Value = 17
$ws.Range("AK8").Value = 34
$ws.Range("AM8").Value = 29
$ws.Range("AO8").Value = 12
$ws.Range("AW8").Value = 5
$ws.Range("AX8").Value = 17

# Row 9
$ws.Range("Z9").Value = 13
$ws.Range("AN9").Value = 3.75

# Row 10
$ws.Range("G10").Value = 3.9
$ws.Range("I10").Value = 1.95
$ws.Range("L10").Value = 2.6
$ws.Range("U10").Value = 1.73
$ws.Range("V10").Value = 2
$ws.Range("AI10").Value = 9.5
$ws.Range("AO10").Value = 21
$ws.Range("BD10").Value = 151

# Row 11
$ws.Range("G11").Value = 2.7
$ws.Range("I11").Value = 2.6
$ws.Range("AZ11").Value = 41

# Row 12
$ws.Range("G12").Value = 1.62
$ws.Range("H12").Value = 4
$ws.Range("I12").Value = 5
$ws.Range("AC12").Value = 10
$ws.Range("AE12").Value = 19
$ws.Range("AH12").Value = 13
$ws.Range("AI12").Value = 26
$ws.Range("AJ12").Value = 17
$ws.Range("AQ12").Value = 26
$ws.Range("AW12").Value = 7
$ws.Range("BC12").Value = 151

# Row 13
$ws.Range("G13").Value = 2.15
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 3.2
$ws.Range("J13").Value = 2.88
$ws.Range("L13").Value = 3.75
$ws.Range("U13").Value = 1.73
$ws.Range("V13").Value = 2
$ws.Range("W13").Value = 8.5
$ws.Range("X13").Value = 11
$ws.Range("Y13").Value = 9
$ws.Range("Z13").Value = 21
$ws.Range("AA13").Value = 17
$ws.Range("AH13").Value = 11
$ws.Range("AI13").Value = 17
$ws.Range("AK13").Value = 34
$ws.Range("AL13").Value = 26
$ws.Range("AM13").Value = 34
$ws.Range("AN13").Value = 4.33
$ws.Range("AO13").Value = 12
$ws.Range("AX13").Value = 17

# Row 14
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("O14").Value = 1.25
$ws.Range("P14").Value = 3.75
$ws.Range("Q14").Value = 1.83
$ws.Range("R14").Value = 2.03

# Row 15
$ws.Range("G15").Value = 1.67
$ws.Range("H15").Value = 3.7
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 2.25
$ws.Range("U15").Value = 1.75
$ws.Range("V15").Value = 2
$ws.Range("X15").Value = 8.5
$ws.Range("AI15").Value = 26
$ws.Range("AJ15").Value = 15

# Row 17
$ws.Range("G17").Value = 3.3
$ws.Range("N17").Value = 8.5
$ws.Range("AI17").Value = 9.5
$ws.Range("AJ17").Value = 9
$ws.Range("AK17").Value = 19
$ws.Range("AR17").Value = 101
$ws.Range("AU17").Value = 8.5
$ws.Range("AX17").Value = 12

# Row 20
$ws.Range("I20").Value = 3.7
$ws.Range("K20").Value = 2.1
$ws.Range("O20").Value = 1.33
$ws.Range("P20").Value = 3.25
$ws.Range("Q20").Value = 2.05
$ws.Range("R20").Value = 1.75
$ws.Range("S20").Value = 1.44
$ws.Range("T20").Value = 2.63
$ws.Range("U20").Value = 1.91
$ws.Range("V20").Value = 1.91
$ws.Range("W20").Value = 7
$ws.Range("X20").Value = 9
$ws.Range("AB20").Value = 29
$ws.Range("AC20").Value = 9.5
$ws.Range("AG20").Value = 301
$ws.Range("AH20").Value = 10
$ws.Range("AM20").Value = 41
$ws.Range("AT20").Value = 2.63
$ws.Range("BA20").Value = 101

# Row 21
$ws.Range("G21").Value = 4.75
$ws.Range("H21").Value = 3.5
$ws.Range("I21").Value = 1.75
$ws.Range("J21").Value = 5
$ws.Range("K21").Value = 2.2
$ws.Range("L21").Value = 2.4
$ws.Range("M21").Value = 1.07
$ws.Range("N21").Value = 9
$ws.Range("U21").Value = 1.95
$ws.Range("V21").Value = 1.8
$ws.Range("W21").Value = 12
$ws.Range("X21").Value = 23
$ws.Range("Z21").Value = 51
$ws.Range("AA21").Value = 41
$ws.Range("AD21").Value = 7
$ws.Range("AE21").Value = 17
$ws.Range("AG21").Value = 351
$ws.Range("AH21").Value = 6.5
$ws.Range("AI21").Value = 8
$ws.Range("AK21").Value = 13
$ws.Range("AN21").Value = 6.5
$ws.Range("AO21").Value = 26
$ws.Range("AQ21").Value = 101
$ws.Range("AR21").Value = 126
$ws.Range("AW21").Value = 3.6
$ws.Range("AX21").Value = 9.5
$ws.Range("AZ21").Value = 29

# Row 23
$ws.Range("I23").Value = 2.7
$ws.Range("X23").Value = 12
$ws.Range("AA23").Value = 26
$ws.Range("AD23").Value = 6
$ws.Range("AU23").Value = 9
$ws.Range("AW23").Value = 4.5

# Row 24
$ws.Range("K24").Value = 2.2
$ws.Range("L24").Value = 3.4
$ws.Range("N24").Value = 13
$ws.Range("O24").Value = 1.25
$ws.Range("P24").Value = 3.75
$ws.Range("Q24").Value = 1.83
$ws.Range("R24").Value = 2.03
$ws.Range("U24").Value = 1.67
$ws.Range("V24").Value = 2.1
$ws.Range("AE24").Value = 13
$ws.Range("AH24").Value = 10
$ws.Range("AJ24").Value = 10
$ws.Range("AM24").Value = 29
$ws.Range("AY24").Value = 23
